$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.669.46'
$ws.Range('E2').Value = '  -2.30%  '

# Row 3
$ws.Range('D3').Value = '3.484.13'
$ws.Range('E3').Value = '  -3.62%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.36%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '580.18'
$ws.Range('E5').Value = '  -3.69%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '189.27'
$ws.Range('E6').Value = '  -3.37%  '

# Row 7
$ws.Range('D7').Value = '3.467.52'
$ws.Range('E7').Value = '  -3.76%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  -3.94%  '

# Row 9
$ws.Range('E9').Value = '  +0.10%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.201'
$ws.Range('E10').Value = '  -4.97%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.613'
$ws.Range('E11').Value = '  -5.35%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '51.58'
$ws.Range('E12').Value = '  -3.17%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000283'
$ws.Range('E13').Value = '  -6.81%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.07'
$ws.Range('E14').Value = '  -5.48%  '

# Row 15
$ws.Range('D15').Value = '4.030.38'
$ws.Range('E15').Value = '  -3.73%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '634.18'

# Row 17
$ws.Range('D17').Value = '68.859.59'
$ws.Range('E17').Value = '  -2.18%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.37'
$ws.Range('E18').Value = '  -4.66%  '

# Row 19
$ws.Range('D19').Value = '3.459.74'
$ws.Range('E19').Value = '  -4.30%  '

# Row 20
$ws.Range('E20').Value = '  -2.36%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '18.09'
$ws.Range('E21').Value = '  -5.12%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.944'
$ws.Range('E22').Value = '  -5.78%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.74'
$ws.Range('E23').Value = '  -4.60%  '

# Row 24
$ws.Range('E24').Value = '  +3.67%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '99.56'
$ws.Range('E25').Value = '  -3.73%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '4.30'
$ws.Range('E26').Value = '  -6.88%  '

# Row 27
$ws.Range('E27').Value = '  -4.68%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.06'
$ws.Range('E28').Value = '  +2.16%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '10.03'
$ws.Range('E29').Value = '  -5.70%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.20'
$ws.Range('E30').Value = '  -5.39%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '32.44'
$ws.Range('E31').Value = '  -4.08%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.69'
$ws.Range('E32').Value = '  -8.22%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.05'
$ws.Range('E33').Value = '  -13.88%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.55'
$ws.Range('E34').Value = '  -6.02%  '

# Row 35
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '60.77'
$ws.Range('E35').Value = '  -4.07%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.108'
$ws.Range('E36').Value = '  -7.77%  '

# Row 37
$ws.Range('D37').Value = '3.706.91'
$ws.Range('E37').Value = '  -5.85%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.32%  '

# Row 39
$ws.Range('D39').Value = '0.0₃0787'
$ws.Range('E39').Value = '  -11.16%  '

# Row 40
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '502.39'
$ws.Range('E40').Value = '  -5.85%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.55'
$ws.Range('E41').Value = '  +0.48%  '

# Row 42
$ws.Range('E42').Value = '  -4.08%  '

# Row 43
$ws.Range('E43').Value = '  -5.63%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.132'
$ws.Range('E44').Value = '  -1.35%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '34.14'
$ws.Range('E45').Value = '  -7.43%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0438'
$ws.Range('E46').Value = '  -4.92%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.35'
$ws.Range('E47').Value = '  -7.02%  '

# Row 48
$ws.Range('E48').Value = '  -2.95%  '

# Row 49
$ws.Range('E49').Value = '  -4.37%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  -0.03%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.05'
$ws.Range('E51').Value = '  -6.02%  '
